$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginTestData")

# Normalize header casing: "Username" -> "username", "Password" -> "password"
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"

# Update the active selection to B1
$ws.Activate()
$ws.Range("B1").Select()
